$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'301.51"
$ws.Range("E2").Value = "'-4.11%"
$ws.Range("G2").Value = "'14"
$ws.Range("D3").Value = "'35.53"
$ws.Range("E3").Value = "'1.59%"
$ws.Range("G3").Value = "'14"
$ws.Range("D4").Value = "'5.053"
$ws.Range("E4").Value = "'-1.26%"
$ws.Range("G4").Value = "'14"
$ws.Range("D5").Value = "'0.07997"
$ws.Range("E5").Value = "'-1.89%"
$ws.Range("G5").Value = "'14"
$ws.Range("D6").Value = "'1.921"
$ws.Range("E6").Value = "'-9.88%"
$ws.Range("G6").Value = "'14"
$ws.Range("D7").Value = "'7.816"
$ws.Range("E7").Value = "'-1.77%"
$ws.Range("G7").Value = "'14"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.9249"
$ws.Range("E8").Value = "'-0.80%"
$ws.Range("G8").Value = "'14"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "'0.1397"
$ws.Range("E9").Value = "'34.18%"
$ws.Range("G9").Value = "'14"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1891"
$ws.Range("E10").Value = "'-1.87%"
$ws.Range("G10").Value = "'14"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.09247"
$ws.Range("E11").Value = "'1.59%"
$ws.Range("G11").Value = "'14"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.03381"
$ws.Range("E12").Value = "'-6.99%"
$ws.Range("G12").Value = "'14"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09872"
$ws.Range("E13").Value = "'-0.14%"
$ws.Range("G13").Value = "'14"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001394"
$ws.Range("E14").Value = "'-2.70%"
$ws.Range("G14").Value = "'14"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.005849"
$ws.Range("E15").Value = "'1.30%"
$ws.Range("G15").Value = "'14"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.505"
$ws.Range("E16").Value = "'1.03%"
$ws.Range("G16").Value = "'14"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "'4.062"
$ws.Range("E17").Value = "'-1.98%"
$ws.Range("G17").Value = "'14"
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D18").Value = "'2.959"
$ws.Range("E18").Value = "'5.20%"
$ws.Range("G18").Value = "'14"
$ws.Range("D19").Value = "'0.3402"
$ws.Range("E19").Value = "'-0.20%"
$ws.Range("G19").Value = "'14"
$ws.Range("G20").Value = "'14"
$ws.Range("D21").Value = "'5.044"
$ws.Range("E21").Value = "'-1.00%"
$ws.Range("G21").Value = "'14"
$ws.Range("D22").Value = "'0.2405"
$ws.Range("E22").Value = "'8.58%"
$ws.Range("G22").Value = "'14"
$ws.Range("D23").Value = "'0.04492"
$ws.Range("E23").Value = "'-1.60%"
$ws.Range("G23").Value = "'14"
$ws.Range("D24").Value = "'0.001215"
$ws.Range("E24").Value = "'-2.65%"
$ws.Range("G24").Value = "'14"
$ws.Range("D25").Value = "'0.004791"
$ws.Range("E25").Value = "'2.11%"
$ws.Range("G25").Value = "'14"
$ws.Range("D26").Value = "'0.0001251"
$ws.Range("E26").Value = "'0.03%"
$ws.Range("G26").Value = "'14"
$ws.Range("D27").Value = "'0.0003006"
$ws.Range("E27").Value = "'-33.31%"
$ws.Range("G27").Value = "'14"
$ws.Range("G28").Value = "'14"
$ws.Range("G29").Value = "'14"
$ws.Range("G30").Value = "'14"
$ws.Range("G31").Value = "'14"
$ws.Range("G32").Value = "'14"
$ws.Range("G33").Value = "'14"
$ws.Range("G34").Value = "'14"
$ws.Range("G35").Value = "'14"
$ws.Range("G36").Value = "'14"
$ws.Range("G37").Value = "'14"
$ws.Range("G38").Value = "'14"
$ws.Range("D39").Value = "'0.01914"
$ws.Range("E39").Value = "'-2.17%"
$ws.Range("G39").Value = "'14"
$ws.Range("D40").Value = "'0.04745"
$ws.Range("E40").Value = "'-3.04%"
$ws.Range("G40").Value = "'14"
$ws.Range("D41").Value = "'0.007355"
$ws.Range("E41").Value = "'-3.38%"
$ws.Range("G41").Value = "'14"
$ws.Range("D42").Value = "'0.009682"
$ws.Range("E42").Value = "'22.89%"
$ws.Range("G42").Value = "'14"
$ws.Range("E43").Value = "'-4.10%"
$ws.Range("G43").Value = "'14"
$ws.Range("D44").Value = "'0.002112"
$ws.Range("E44").Value = "'0.39%"
$ws.Range("G44").Value = "'14"
$ws.Range("D45").Value = "'0.01055"
$ws.Range("E45").Value = "'-10.16%"
$ws.Range("G45").Value = "'14"
$ws.Range("D46").Value = "'0.00006264"
$ws.Range("E46").Value = "'-7.27%"
$ws.Range("G46").Value = "'14"
$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("E47").Value = "'0.14%"
$ws.Range("G47").Value = "'14"
$ws.Range("E48").Value = "'-65.05%"
$ws.Range("G48").Value = "'14"
$ws.Range("E49").Value = "'-2.39%"
$ws.Range("G49").Value = "'14"
$ws.Range("D50").Value = "'0.00002105"
$ws.Range("E50").Value = "'0.14%"
$ws.Range("G50").Value = "'14"
$ws.Range("E51").Value = "'0.14%"
$ws.Range("G51").Value = "'14"

Write-Host "Applied 144 cell updates"
